$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# re-run RU 1001; without crop
$ws.Range("B2").Value = 55.3352948193206
$ws.Range("L2").Value = 48.6313846276599
$ws.Range("B3").Value = 40.563158268296
$ws.Range("L3").Value = 36.9618773883661
